$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(18, 8).Value2 = 8118.6
$ws.Cells.Item(18, 9).Value2 = 3616
$ws.Cells.Item(18, 11).Value2 = 3616
$ws.Cells.Item(18, 13).Value2 = -3332

$ws.Cells.Item(43, 8).Value2 = 916.04
$ws.Cells.Item(43, 9).Value2 = 718.5
$ws.Cells.Item(43, 10).Value2 = 1009
$ws.Cells.Item(43, 11).Value2 = 718.5
$ws.Cells.Item(43, 12).Value2 = 1009
$ws.Cells.Item(43, 13).Value2 = -649.5
$ws.Cells.Item(43, 14).Value2 = -1147

$ws.Cells.Item(74, 8).Value2 = 3000
$ws.Cells.Item(74, 9).Value2 = 3000
$ws.Cells.Item(74, 10).Value2 = 0
$ws.Cells.Item(74, 11).Value2 = 3000
$ws.Cells.Item(74, 12).Value2 = 0
$ws.Cells.Item(74, 13).ClearContents()
$ws.Cells.Item(74, 14).Value2 = -2064

$ws.Cells.Item(77, 8).Value2 = 3000
$ws.Cells.Item(77, 9).Value2 = 3000
$ws.Cells.Item(77, 10).Value2 = 0
$ws.Cells.Item(77, 11).Value2 = 15000
$ws.Cells.Item(77, 12).Value2 = 0
$ws.Cells.Item(77, 13).ClearContents()
$ws.Cells.Item(77, 14).Value2 = -10320

$ws.Cells.Item(106, 8).Value2 = 2086
$ws.Cells.Item(106, 9).Value2 = 2086
$ws.Cells.Item(106, 11).Value2 = 2086
$ws.Cells.Item(106, 13).Value2 = -1455

$ws.Cells.Item(121, 8).Value2 = 1322.3334
$ws.Cells.Item(121, 10).Value2 = 1556.8
$ws.Cells.Item(121, 12).Value2 = 4670.4
$ws.Cells.Item(121, 14).Value2 = -8164.4

$ws.Cells.Item(137, 8).Value2 = 1069.9744
$ws.Cells.Item(137, 9).Value2 = 814.26666
$ws.Cells.Item(137, 11).Value2 = 2442.79998
$ws.Cells.Item(137, 13).Value2 = 107.2000200000002

$ws.Cells.Item(138, 8).Value2 = 2049.9524
$ws.Cells.Item(138, 9).Value2 = 2071.4412
$ws.Cells.Item(138, 10).Value2 = 1958.625
$ws.Cells.Item(138, 11).Value2 = 6214.323600000001
$ws.Cells.Item(138, 12).Value2 = 5875.875
$ws.Cells.Item(138, 13).Value2 = -1074.323600000001
$ws.Cells.Item(138, 14).Value2 = -16155.875

$ws.Cells.Item(139, 8).Value2 = 50462.375
$ws.Cells.Item(139, 10).Value2 = 50462.375
$ws.Cells.Item(139, 12).Value2 = 50462.375
$ws.Cells.Item(139, 14).Value2 = -60742.375

$ws.Cells.Item(140, 8).Value2 = 57240.4
$ws.Cells.Item(140, 10).Value2 = 57240.4
$ws.Cells.Item(140, 12).Value2 = 57240.4
$ws.Cells.Item(140, 14).Value2 = -67600.39999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value2 = 2666.3586
$ws.Cells.Item(32, 9).Value2 = 2031.4938
$ws.Cells.Item(32, 10).Value2 = 7341.273
$ws.Cells.Item(32, 11).Value2 = 2031.4938
$ws.Cells.Item(32, 12).Value2 = 7341.273
$ws.Cells.Item(32, 13).Value2 = -1744.4938
$ws.Cells.Item(32, 14).Value2 = -7915.273

$ws.Cells.Item(61, 8).Value2 = 2819.9565
$ws.Cells.Item(61, 9).Value2 = 2293
$ws.Cells.Item(61, 11).Value2 = 2293
$ws.Cells.Item(61, 13).Value2 = -2081

$ws.Cells.Item(74, 8).Value2 = 1064.9535
$ws.Cells.Item(74, 9).Value2 = 822.5454999999999
$ws.Cells.Item(74, 10).Value2 = 1864.9
$ws.Cells.Item(74, 11).Value2 = 822.5454999999999
$ws.Cells.Item(74, 12).Value2 = 1864.9
$ws.Cells.Item(74, 13).Value2 = 51.45450000000005
$ws.Cells.Item(74, 14).Value2 = -3612.9

$ws.Cells.Item(77, 8).Value2 = 1064.9535
$ws.Cells.Item(77, 9).Value2 = 822.5454999999999
$ws.Cells.Item(77, 10).Value2 = 1864.9
$ws.Cells.Item(77, 11).Value2 = 4112.7275
$ws.Cells.Item(77, 12).Value2 = 9324.5
$ws.Cells.Item(77, 13).Value2 = 255.2725
$ws.Cells.Item(77, 14).Value2 = -18060.5

$ws.Cells.Item(122, 8).Value2 = 2054.8572
$ws.Cells.Item(122, 9).Value2 = 1685.875
$ws.Cells.Item(122, 10).Value2 = 2546.8333
$ws.Cells.Item(122, 11).Value2 = 5057.625
$ws.Cells.Item(122, 12).Value2 = 7640.499899999999
$ws.Cells.Item(122, 13).Value2 = -2607.625
$ws.Cells.Item(122, 14).Value2 = -12540.4999

$ws.Cells.Item(132, 8).Value2 = 1461.4807
$ws.Cells.Item(132, 9).Value2 = 1034.7354
$ws.Cells.Item(132, 10).Value2 = 2267.5557
$ws.Cells.Item(132, 11).Value2 = 3104.2062
$ws.Cells.Item(132, 12).Value2 = 6802.6671
$ws.Cells.Item(132, 13).Value2 = -574.2062000000001
$ws.Cells.Item(132, 14).Value2 = -11862.6671

$ws.Cells.Item(136, 8).Value2 = 2819.9565
$ws.Cells.Item(136, 9).Value2 = 2293
$ws.Cells.Item(136, 11).Value2 = 6879
$ws.Cells.Item(136, 13).Value2 = -4329

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value2 = 4512.375
$ws.Cells.Item(20, 9).Value2 = 4416.6665
$ws.Cells.Item(20, 11).Value2 = 4416.6665
$ws.Cells.Item(20, 13).Value2 = -4169.6665

$ws.Cells.Item(99, 8).Value2 = 3449
$ws.Cells.Item(99, 9).Value2 = 3498
$ws.Cells.Item(99, 10).Value2 = 3432.6667
$ws.Cells.Item(99, 11).Value2 = 3498
$ws.Cells.Item(99, 12).Value2 = 3432.6667
$ws.Cells.Item(99, 13).Value2 = -2000
$ws.Cells.Item(99, 14).Value2 = -6428.6667

$ws.Cells.Item(107, 8).Value2 = 1516.5
$ws.Cells.Item(107, 9).Value2 = 1459.8
$ws.Cells.Item(107, 10).Value2 = 1800
$ws.Cells.Item(107, 11).Value2 = 1459.8
$ws.Cells.Item(107, 12).Value2 = 1800
$ws.Cells.Item(107, 13).Value2 = 460.2
$ws.Cells.Item(107, 14).Value2 = -5640

$ws.Cells.Item(134, 8).Value2 = 9254.583000000001
$ws.Cells.Item(134, 9).Value2 = 11214.134
$ws.Cells.Item(134, 11).Value2 = 33642.402
$ws.Cells.Item(134, 13).Value2 = -31107.402

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(58, 8).Value2 = 1146086.9
$ws.Cells.Item(58, 9).Value2 = 1977631.2
$ws.Cells.Item(58, 10).Value2 = 2713.25
$ws.Cells.Item(58, 11).Value2 = 1977631.2
$ws.Cells.Item(58, 12).Value2 = 2713.25
$ws.Cells.Item(58, 13).Value2 = -1977428.2
$ws.Cells.Item(58, 14).Value2 = -3119.25

$ws.Cells.Item(136, 8).Value2 = 1146086.9
$ws.Cells.Item(136, 9).Value2 = 1977631.2
$ws.Cells.Item(136, 10).Value2 = 2713.25
$ws.Cells.Item(136, 11).Value2 = 5932893.6
$ws.Cells.Item(136, 12).Value2 = 8139.75
$ws.Cells.Item(136, 13).Value2 = -5930343.6
$ws.Cells.Item(136, 14).Value2 = -13239.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(11, 8).Value2 = 624
$ws.Cells.Item(11, 10).Value2 = 110
$ws.Cells.Item(11, 12).Value2 = 330
$ws.Cells.Item(11, 14).Value2 = -610

$ws.Cells.Item(23, 8).Value2 = 200
$ws.Cells.Item(23, 10).Value2 = 250
$ws.Cells.Item(23, 12).Value2 = 750
$ws.Cells.Item(23, 14).Value2 = -1220

$ws.Cells.Item(40, 8).Value2 = 512.5
$ws.Cells.Item(40, 9).Value2 = 500
$ws.Cells.Item(40, 10).Value2 = 550
$ws.Cells.Item(40, 11).Value2 = 2000
$ws.Cells.Item(40, 12).Value2 = 2200
$ws.Cells.Item(40, 13).Value2 = -1931
$ws.Cells.Item(40, 14).Value2 = -2338

$ws.Cells.Item(131, 8).Value2 = 15269.109
$ws.Cells.Item(131, 9).Value2 = 352.66666
$ws.Cells.Item(131, 10).Value2 = 19431.838
$ws.Cells.Item(131, 11).Value2 = 1057.99998
$ws.Cells.Item(131, 12).Value2 = 58295.514
$ws.Cells.Item(131, 13).Value2 = 3982.00002
$ws.Cells.Item(131, 14).Value2 = -68375.514

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value2 = 19319.4
$ws.Cells.Item(70, 9).Value2 = 29332.334
$ws.Cells.Item(70, 11).Value2 = 29332.334
$ws.Cells.Item(70, 13).Value2 = -29062.334

$ws.Cells.Item(73, 8).Value2 = 19319.4
$ws.Cells.Item(73, 9).Value2 = 29332.334
$ws.Cells.Item(73, 11).Value2 = 29332.334
$ws.Cells.Item(73, 13).Value2 = -28396.334

$ws.Cells.Item(113, 8).Value2 = 1169
$ws.Cells.Item(113, 9).Value2 = 995.5
$ws.Cells.Item(113, 10).Value2 = 1255.75
$ws.Cells.Item(113, 11).Value2 = 995.5
$ws.Cells.Item(113, 12).Value2 = 1255.75
$ws.Cells.Item(113, 13).Value2 = 1174.5
$ws.Cells.Item(113, 14).Value2 = -5595.75

$ws.Cells.Item(137, 8).Value2 = 61786.668
$ws.Cells.Item(137, 10).Value2 = 61786.668
$ws.Cells.Item(137, 12).Value2 = 61786.668
$ws.Cells.Item(137, 14).Value2 = -71986.66800000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(22, 8).Value2 = 3265.4285
$ws.Cells.Item(22, 9).Value2 = 3239.5
$ws.Cells.Item(22, 10).Value2 = 3300
$ws.Cells.Item(22, 11).Value2 = 3239.5
$ws.Cells.Item(22, 12).Value2 = 3300
$ws.Cells.Item(22, 13).Value2 = -2944.5
$ws.Cells.Item(22, 14).Value2 = -3890

$ws.Cells.Item(27, 8).Value2 = 3265.4285
$ws.Cells.Item(27, 9).Value2 = 3239.5
$ws.Cells.Item(27, 10).Value2 = 3300
$ws.Cells.Item(27, 11).Value2 = 3239.5
$ws.Cells.Item(27, 12).Value2 = 3300
$ws.Cells.Item(27, 13).Value2 = -3132.5
$ws.Cells.Item(27, 14).Value2 = -3514

$ws.Cells.Item(40, 8).Value2 = 6839.08
$ws.Cells.Item(40, 9).Value2 = 7269.647
$ws.Cells.Item(40, 10).Value2 = 5924.125
$ws.Cells.Item(40, 11).Value2 = 7269.647
$ws.Cells.Item(40, 12).Value2 = 5924.125
$ws.Cells.Item(40, 13).Value2 = -7133.647
$ws.Cells.Item(40, 14).Value2 = -6196.125

$ws.Cells.Item(46, 8).Value2 = 1539.8572
$ws.Cells.Item(46, 9).Value2 = 1000.5
$ws.Cells.Item(46, 10).Value2 = 1755.6
$ws.Cells.Item(46, 11).Value2 = 1000.5
$ws.Cells.Item(46, 12).Value2 = 1755.6
$ws.Cells.Item(46, 13).Value2 = -812.5
$ws.Cells.Item(46, 14).Value2 = -2131.6

$ws.Cells.Item(122, 8).Value2 = 4998.846
$ws.Cells.Item(122, 9).Value2 = 6835.2
$ws.Cells.Item(122, 11).Value2 = 20505.6
$ws.Cells.Item(122, 13).Value2 = -18055.6

$ws.Cells.Item(132, 8).Value2 = 1410.7778
$ws.Cells.Item(132, 9).Value2 = 957.11365
$ws.Cells.Item(132, 11).Value2 = 2871.34095
$ws.Cells.Item(132, 13).Value2 = -341.3409499999998

$ws.Cells.Item(136, 8).Value2 = 2458.3713
$ws.Cells.Item(136, 9).Value2 = 1841.4348
$ws.Cells.Item(136, 10).Value2 = 3640.8333
$ws.Cells.Item(136, 11).Value2 = 5524.3044
$ws.Cells.Item(136, 12).Value2 = 10922.4999
$ws.Cells.Item(136, 13).Value2 = -2974.3044
$ws.Cells.Item(136, 14).Value2 = -16022.4999

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(98, 8).Value2 = 0
$ws.Cells.Item(98, 10).Value2 = 0
$ws.Cells.Item(98, 12).ClearContents()
$ws.Cells.Item(98, 14).ClearContents()

$ws.Cells.Item(122, 8).Value2 = 78839.7
$ws.Cells.Item(122, 9).Value2 = 87504.89
$ws.Cells.Item(122, 11).Value2 = 262514.67
$ws.Cells.Item(122, 13).Value2 = -260064.67

$ws.Cells.Item(126, 8).Value2 = 7331
$ws.Cells.Item(126, 9).Value2 = 8430.5
$ws.Cells.Item(126, 11).Value2 = 25291.5
$ws.Cells.Item(126, 13).Value2 = -22821.5

$ws.Cells.Item(132, 8).Value2 = 1751.6888
$ws.Cells.Item(132, 9).Value2 = 1215.04
$ws.Cells.Item(132, 10).Value2 = 2422.5
$ws.Cells.Item(132, 11).Value2 = 3645.12
$ws.Cells.Item(132, 12).Value2 = 7267.5
$ws.Cells.Item(132, 13).Value2 = -1115.12
$ws.Cells.Item(132, 14).Value2 = -12327.5

$ws.Cells.Item(136, 8).Value2 = 12628038
$ws.Cells.Item(136, 9).Value2 = 19843172
$ws.Cells.Item(136, 11).Value2 = 59529516
$ws.Cells.Item(136, 13).Value2 = -59526966
